$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove E4 (old "0010?" string); no longer used after this edit ---
$ws.Range("E4").Clear()

$ws.Range("I4").Value = 0
$ws.Range("D5").Value = 1
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 1
$ws.Range("D6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 1
$ws.Range("I6").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = 1
$ws.Range("D8").Value = "1 [use the appropriately extended immediate as source B]"
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 1
$ws.Range("I8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("D12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 0
$ws.Range("D13").Value = 1
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 1
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 1
$ws.Range("D15").Value = 1
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("D16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 1
$ws.Range("I16").Value = 0
$ws.Range("D17").Value = 1
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = 1
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 1
$ws.Range("I18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 0
$ws.Range("D21").Value = 1
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = "-"
$ws.Range("D22").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 0
$ws.Range("D23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 1
$ws.Range("I23").Value = 0
$ws.Range("D24").Value = 1
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = "-"
$ws.Range("D25").Value = 1
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = "-"
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = "-"
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = "-"
$ws.Range("D28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 0
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = "-"

# --- Update the active selection in the bottom-right frozen pane to G20 ---
$ws.Range("G20").Select()

Write-Output "edit applied"
